# Arreglado bugs importar siempre la misma rejilla
# Remove accented characters from a handful of text labels in the grid so
# the import routine stops choking on them.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("K1").Value = "Angela"
$ws.Range("M1").Value = "Rocio"
$ws.Range("O13").Value = "Energico"
$ws.Range("A14").Value = "Egocentrico"
$ws.Range("O14").Value = "Empatico"
$ws.Range("A17").Value = "Frio"
$ws.Range("O17").Value = "Calido"
$ws.Range("A19").Value = "Antipatico"
$ws.Range("O20").Value = "Grunyon"
